# "add genders in batumi" - append a new year column (2023) of trade
# statistics to the Imereti region trade worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Write the new 2023 values into column S -------------------------------
$ws.Range("S3").Value = 2023

$ws.Range("S4").Value = 2676.3
$ws.Range("S5").Value = 451.7
$ws.Range("S6").Value = 15011
$ws.Range("S7").Value = 10515
$ws.Range("S8").Value = 764.7
$ws.Range("S9").Value = 85.6
$ws.Range("S10").Value = 97.9
$ws.Range("S11").Value = 366.1
$ws.Range("S12").Value = 10.9
$ws.Range("S13").Value = 2344.6
$ws.Range("S14").Value = 2274.9

# --- Carry over the same formatting column R uses onto column S ------------
# (header band, year-row style, and the three numeric bands used by the
# other columns) so the new column looks like the rest of the table.
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("R3:R14").Copy()
$ws.Range("S3:S14").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = $false

# --- Extend the title merge from A1:R1 to A1:S1 -----------------------------
$ws.Range("A1:R1").UnMerge()
$ws.Range("A1:S1").Merge()

# --- Update the view so the new column is the active selection -------------
# Scroll the viewport so column I is the first visible column, then leave
# the final selection on S3:S14 (matching the saved sheetView/selection).
$ws.Range("I1").Select()
try { $excel.ActiveWindow.ScrollColumn = 9 } catch {}
$ws.Range("S3:S14").Select()

# --- Resize the saved window (cosmetic; best effort) ------------------------
try {
    $excel.ActiveWindow.Width = 22215
    $excel.ActiveWindow.Height = 9285
} catch {}
